$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 699.25
$ws.Range("I2").Value = 337.8
$ws.Range("K2").Value = 337.8
$ws.Range("M2").Value = -224.8
$ws.Range("H33").Value = 220.6
$ws.Range("I33").Value = 195.45454
$ws.Range("J33").Value = 251.33333
$ws.Range("K33").Value = 195.45454
$ws.Range("L33").Value = 251.33333
$ws.Range("M33").Value = 33.54545999999999
$ws.Range("N33").Value = -709.3333299999999
$ws.Range("H40").Value = 3993.0625
$ws.Range("I40").Value = 1974.75
$ws.Range("K40").Value = 1974.75
$ws.Range("M40").Value = -1799.75
$ws.Range("H129").Value = 58824604
$ws.Range("I129").Value = 83333816
$ws.Range("K129").Value = 250001448
$ws.Range("M129").Value = -249996448
$ws.Range("H132").Value = 10029.417
$ws.Range("I132").Value = 3034.9333
$ws.Range("J132").Value = 45001.832
$ws.Range("K132").Value = 9104.7999
$ws.Range("L132").Value = 135005.496
$ws.Range("M132").Value = -6574.7999
$ws.Range("N132").Value = -140065.496
$ws.Range("H137").Value = 2127.484
$ws.Range("I137").Value = 2076.9285
$ws.Range("K137").Value = 6230.7855
$ws.Range("M137").Value = -3680.7855
$ws.Range("H138").Value = 394465.3
$ws.Range("J138").Value = 716287.8
$ws.Range("L138").Value = 2148863.4
$ws.Range("N138").Value = -2159143.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 256.42856
$ws.Range("I22").Value = 256.42856
$ws.Range("K22").Value = 256.42856
$ws.Range("M22").Value = 42.57144
$ws.Range("H32").Value = 7236.2446
$ws.Range("I32").Value = 7396.6904
$ws.Range("K32").Value = 7396.6904
$ws.Range("M32").Value = -7109.6904
$ws.Range("H45").Value = 2252.8696
$ws.Range("I45").Value = 1987.5
$ws.Range("J45").Value = 2665.6667
$ws.Range("K45").Value = 1987.5
$ws.Range("L45").Value = 2665.6667
$ws.Range("M45").Value = -1610.5
$ws.Range("N45").Value = -3419.6667
$ws.Range("H61").Value = 5818.796
$ws.Range("I61").Value = 3711.1614
$ws.Range("K61").Value = 3711.1614
$ws.Range("M61").Value = -3499.1614
$ws.Range("H74").Value = 2735.9211
$ws.Range("I74").Value = 2217.394
$ws.Range("J74").Value = 6158.2
$ws.Range("K74").Value = 2217.394
$ws.Range("L74").Value = 6158.2
$ws.Range("M74").Value = -1343.394
$ws.Range("N74").Value = -7906.2
$ws.Range("H77").Value = 2735.9211
$ws.Range("I77").Value = 2217.394
$ws.Range("J77").Value = 6158.2
$ws.Range("K77").Value = 11086.97
$ws.Range("L77").Value = 30791
$ws.Range("M77").Value = -6718.969999999999
$ws.Range("N77").Value = -39527
$ws.Range("H102").Value = 4247.1333
$ws.Range("J102").Value = 7227.1665
$ws.Range("L102").Value = 7227.1665
$ws.Range("N102").Value = -10471.1665
$ws.Range("H136").Value = 5818.796
$ws.Range("I136").Value = 3711.1614
$ws.Range("K136").Value = 11133.4842
$ws.Range("M136").Value = -8583.484199999999
$ws.Range("H137").Value = 90000
$ws.Range("J137").Value = 97142.86
$ws.Range("L137").Value = 97142.86
$ws.Range("N137").Value = -107342.86

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 2667.7144
$ws.Range("J37").Value = 4407
$ws.Range("L37").Value = 4407
$ws.Range("N37").Value = -4681
$ws.Range("H105").Value = 3368
$ws.Range("I105").Value = 2335.9048
$ws.Range("K105").Value = 2335.9048
$ws.Range("M105").Value = -588.9047999999998

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2247.3
$ws.Range("I31").Value = 1426.2646
$ws.Range("J31").Value = 6899.8335
$ws.Range("K31").Value = 1426.2646
$ws.Range("L31").Value = 6899.8335
$ws.Range("M31").Value = -1131.2646
$ws.Range("N31").Value = -7489.8335
$ws.Range("H34").Value = 2247.3
$ws.Range("I34").Value = 1426.2646
$ws.Range("J34").Value = 6899.8335
$ws.Range("K34").Value = 1426.2646
$ws.Range("L34").Value = 6899.8335
$ws.Range("M34").Value = -1224.2646
$ws.Range("N34").Value = -7303.8335
$ws.Range("H58").Value = 2426
$ws.Range("I58").Value = 2854.8
$ws.Range("K58").Value = 2854.8
$ws.Range("M58").Value = -2651.8
$ws.Range("H105").Value = 2354.4
$ws.Range("I105").Value = 2354.4
$ws.Range("K105").Value = 2354.4
$ws.Range("M105").Value = -607.4000000000001
$ws.Range("H134").Value = 3012.1035
$ws.Range("I134").Value = 1184.725
$ws.Range("J134").Value = 7072.9443
$ws.Range("K134").Value = 3554.175
$ws.Range("L134").Value = 21218.8329
$ws.Range("M134").Value = -1019.175
$ws.Range("N134").Value = -26288.8329
$ws.Range("H136").Value = 2426
$ws.Range("I136").Value = 2854.8
$ws.Range("K136").Value = 8564.400000000001
$ws.Range("M136").Value = -6014.400000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1310
$ws.Range("I3").Value = 1310
$ws.Range("K3").Value = 3930
$ws.Range("M3").Value = -3818
$ws.Range("H5").Value = 6474.421
$ws.Range("J5").Value = 10797.546
$ws.Range("L5").Value = 32392.638
$ws.Range("N5").Value = -32616.638
$ws.Range("H131").Value = 1733.95
$ws.Range("I131").Value = 874.75
$ws.Range("J131").Value = 1948.75
$ws.Range("K131").Value = 2624.25
$ws.Range("L131").Value = 5846.25
$ws.Range("M131").Value = 2415.75
$ws.Range("N131").Value = -15926.25
$ws.Range("H134").Value = 5134.857
$ws.Range("I134").Value = 1486.25
$ws.Range("K134").Value = 4458.75
$ws.Range("M134").Value = 611.25
$ws.Range("H135").Value = 6474.421
$ws.Range("J135").Value = 10797.546
$ws.Range("L135").Value = 97177.914
$ws.Range("N135").Value = -102247.914

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 22999.5
$ws.Range("J44").Value = 22999.5
$ws.Range("L44").Value = 22999.5
$ws.Range("N44").Value = -24191.5
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""
$ws.Range("H122").Value = 4316.143
$ws.Range("I122").Value = 3958.5557
$ws.Range("K122").Value = 11875.6671
$ws.Range("M122").Value = -9425.667099999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4671.35
$ws.Range("J7").Value = 4596.067
$ws.Range("L7").Value = 4596.067
$ws.Range("N7").Value = -4820.067
$ws.Range("H40").Value = 4548.609
$ws.Range("I40").Value = 4526.2354
$ws.Range("J40").Value = 4612
$ws.Range("K40").Value = 4526.2354
$ws.Range("L40").Value = 4612
$ws.Range("M40").Value = -4390.2354
$ws.Range("N40").Value = -4884
$ws.Range("H68").Value = 5483.3335
$ws.Range("I68").Value = 5349.75
$ws.Range("K68").Value = 5349.75
$ws.Range("M68").Value = -4600.75
$ws.Range("H71").Value = 5483.3335
$ws.Range("I71").Value = 5349.75
$ws.Range("K71").Value = 26748.75
$ws.Range("M71").Value = -23004.75
$ws.Range("H126").Value = 4671.35
$ws.Range("J126").Value = 4596.067
$ws.Range("L126").Value = 13788.201
$ws.Range("N126").Value = -18728.201
$ws.Range("H130").Value = 94994
$ws.Range("J130").Value = 94994
$ws.Range("L130").Value = 94994
$ws.Range("N130").Value = -105034
$ws.Range("H136").Value = 4933.579
$ws.Range("I136").Value = 4842.231
$ws.Range("K136").Value = 14526.693
$ws.Range("M136").Value = -11976.693

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 198713.75
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 198713.75
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H81").Value = 47279.707
$ws.Range("J81").Value = 8665
$ws.Range("L81").Value = 17330
$ws.Range("N81").Value = -19452
$ws.Range("H84").Value = 47279.707
$ws.Range("J84").Value = 8665
$ws.Range("L84").Value = 86650
$ws.Range("N84").Value = -97258
$ws.Range("H126").Value = 3641.4285
$ws.Range("I126").Value = 2910.7144
$ws.Range("K126").Value = 8732.143199999999
$ws.Range("M126").Value = -6262.143199999999
$ws.Range("H127").Value = 90000
$ws.Range("J127").Value = 90000
$ws.Range("L127").Value = 90000
$ws.Range("N127").Value = -99920
$ws.Range("H136").Value = 10404.091
$ws.Range("I136").Value = 19817.455
$ws.Range("K136").Value = 59452.36500000001
$ws.Range("M136").Value = -56902.36500000001
